$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in AC1:AE1, matching the
# existing header formatting (bold, centered, bordered) by copying the
# format from the adjacent "Unnamed: 27" header cell (AB1).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Populate the season-record columns for every player row (2-38) with the
# team's season record: 55 wins, 60 losses, 0 ties.
$ws.Range("AC2:AC38").Value = 55
$ws.Range("AD2:AD38").Value = 60
$ws.Range("AE2:AE38").Value = 0
